$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (System Security / No Gaps): tweak wording of the recommendation text ---
$ws.Range("F3").Value = " H&F council must continue to maintain regular testing of its systems to mitigate potential cybersecurity risks. However, any details about the security workers shouldn't be disclosed for security reasons (structure, personnel), as they manage personal and sensitive data."

# --- Row 4 (Access Control / Identity And Access Control): login now requires "correct" email; risk text re-worded ---
$ws.Range("B4").Value = "To log into your personal portal, you have to enter your correct email and password of your account."
$ws.Range("E4").Value = "High – Single-factor login could result in a user's personal data being breached and gathered. Can result in massive consequences."
$ws.Range("E4").Value = "High – Single-factor login could result in a user" + [char]0x2019 + "s personal data being breached and gathered. Can result in massive consequences."

# --- Row 5 (Access Control / Data Security): risk level raised from Low to Medium ---
$ws.Range("E5").Value = "Medium - An unauthorised user can access the device with the active form and tamper with the details entered on the form."

# --- Row 7 (Cookies / Asset Management): Control Area re-labelled ---
$ws.Range("A7").Value = "Data Management"
$ws.Range("F7").Value = "H&F must continue this good practice at all times."

# --- Row 11 (CCTV / Physical Controls): Control Area was a stray header value, now a real label ---
$ws.Range("A11").Value = "Surveilance"

# --- Rows 12-13 (placeholder rows): Control Area re-labelled ---
$ws.Range("A12").Value = "Organisational Procedures"
$ws.Range("A13").Value = "Organisational Procedures"

# --- Adjust the view: scroll the window so D6 is the top-left visible cell, and select H7 ---
$ws.Activate()
$ws.Range("H7").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 4
